$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new data row right above current row 91.
# This pushes the existing rows 91..200 down to 92..201 (and keeps
# all of their values / formatting untouched, matching the diff which
# shows every row from 91 down getting the prior row's data).
$ws.Range("A91:R91").EntireRow.Insert()

# Populate the newly inserted row 91 with the new record's data.
$ws.Range("A91").Value = 11
$ws.Range("B91").Value = "Vega Monumental Concepción"
$ws.Range("C91").Value = "Bíobío"
$ws.Range("D91").Value = 44792
$ws.Range("E91").Value = 8
$ws.Range("F91").Value = 100112003
$ws.Range("G91").Value = "Ajo"
$ws.Range("H91").Value = "Chino"
$ws.Range("I91").Value = "Primera"
$ws.Range("J91").Value = 240
$ws.Range("K91").Value = 21000
$ws.Range("L91").Value = 22000
$ws.Range("M91").Value = 21500
$ws.Range("N91").Value = "$/caja 10 kilos"
$ws.Range("O91").Value = "China"
$ws.Range("P91").Value = 2150
$ws.Range("Q91").Value = 10
$ws.Range("R91").Value = "Hortaliza"
